$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 864
$ws.Range("F3").Value = 13826
$ws.Range("F4").Value = 13608
$ws.Range("F5").Value = 1054
$ws.Range("F7").Value = 42
$ws.Range("F8").Value = 601
$ws.Range("F11").Value = 58
$ws.Range("F12").Value = 766
$ws.Range("F13").Value = 2150
$ws.Range("F14").Value = 110
$ws.Range("F16").Value = 77
$ws.Range("F17").Value = 126
$ws.Range("F21").Value = 409
$ws.Range("F23").Value = 267
$ws.Range("F24").Value = 838
$ws.Range("F25").Value = 95
$ws.Range("F26").Value = 7
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 167
$ws.Range("F7").Value = 1529
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F3").Value = 112
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 864
$ws.Range("F4").Value = 13826
$ws.Range("F5").Value = 13608
$ws.Range("F6").Value = 1054
$ws.Range("F8").Value = 42
$ws.Range("F9").Value = 601
$ws.Range("F12").Value = 58
$ws.Range("F13").Value = 766
$ws.Range("F16").Value = 2150
$ws.Range("F17").Value = 110
$ws.Range("F19").Value = 77
$ws.Range("F20").Value = 126
$ws.Range("F24").Value = 112
$ws.Range("F25").Value = 112
$ws.Range("F28").Value = 409
$ws.Range("F30").Value = 267
$ws.Range("F31").Value = 838
$ws.Range("F32").Value = 167
$ws.Range("F33").Value = 1529
$ws.Range("F37").Value = 95
$ws.Range("F40").Value = 7
